$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D8").Formula = 99
Write-Host "done"
